$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.278.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.509.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.78"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +5.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.53"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.503.14"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.196"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.75%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.579"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.36"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000280"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.080.80"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "617.85"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -9.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.37"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.508.96"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.243.01"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.25"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.25"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.879"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.85"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.29"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.99%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.62"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.24"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.38"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.47"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.11"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.33"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.91"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "570.91"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.59"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Cosmos"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.77"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.22"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.26%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.998"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.406.25"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.326"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0713"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "32.82"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.58"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.84"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.05"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.30%  "
